$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'326.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.66%"
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'43.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.65%"
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'5.548"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.02%"
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'0.08008"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.07%"
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'1.928"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.48%"
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'2.569"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-6.89%"
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.9424"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.55%"
$ws.Range("E8").Style = "Normal"

$ws.Range("E9").Value = "'-3.16%"
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.1843"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-2.94%"
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'11.94"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'40.17%"
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.09593"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.20%"
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'0.04747"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'16.56%"
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.1069"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.28%"
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.001277"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.56%"
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.04074"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-6.39%"
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'0.005911"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.77%"
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'3.377"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-5.66%"
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'4.313"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.29%"
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'-0.33%"
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'0.1403"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'2.89%"
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.2516"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-2.84%"
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'0.001247"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.69%"
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'-7.44%"
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'-2.94%"
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.0003759"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-5.87%"
$ws.Range("E26").Style = "Normal"

$ws.Range("D38").Value = "'0.02541"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-4.82%"
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.05432"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-2.61%"
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'0.007543"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-1.96%"
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.1387"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.59%"
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.007393"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-34.78%"
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.002024"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-4.08%"
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'0.008325"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-11.72%"
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.00007137"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.77%"
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.00000000753"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'0.29%"
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = "'1.29%"
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'0.004814"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'39.65%"
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.00002109"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.29%"
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'0.0002008"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.29%"
$ws.Range("E50").Style = "Normal"
